$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in Week 2 (row 3) timesheet entries ---
# Shared strings get appended in the order the values are first written,
# so write B3, then D3, then C3 to match indices 46/47/48.
$ws.Range("B3").Value = "Figuring out the connection string and conferming it works. Had problems with week 1 tutorial work caused by moving work from PC to laptop and back again."
$ws.Range("D3").Value = "Team meeting talked about problems with database migrations not working."
$ws.Range("C3").Value = "Researched Tag Helpers"

# D3 becomes a highlighted "team meeting" cell like B2/G2/J2/L2/M2 -
# copy that formatting across (reuses the existing highlighted style).
$ws.Range("B2").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# C3 loses its border/fill and just keeps vertical-centered text (no wrap).
$ws.Range("C3").Borders.LineStyle = 0
$ws.Range("C3").WrapText = $false
$ws.Range("C3").VerticalAlignment = -4108

# --- Update the active selection on the sheet ---
$ws.Range("E3").Select() | Out-Null
